$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.745.06'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.603.05'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'211.90"
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').Value = "'19.71"
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').Value = "'0.0847"
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').Value = '1.828.54'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '1.641.39'
$ws.Range('E13').Value = '  +2.88%  '
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = "'209.93"
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = "'7.15"
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  -4.67%  '
$ws.Range('D23').Value = "'9.09"
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = "'143.71"
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = "'1.01"
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = "'0.0508"
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = "'1.16"
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').Value = '1.291.09'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +1.27%  '
$ws.Range('D35').Value = "'1.24"
$ws.Range('E35').Value = '  +20.65%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').Value = "'0.593"
$ws.Range('E37').Value = '  -4.20%  '
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = "'0.828"
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = "'62.98"
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').Value = '1.740.14'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('D45').Value = "'90.49"
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = "'0.102"
$ws.Range('E47').Value = '  +1.16%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = "'6.08"
$ws.Range('E48').Value = '  +18.59%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.0514"
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.53"
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  +0.20%  '
